# Revisions to blueprint and webbrowsing
# Remove the review comment thread (root comment + its reply) that was left
# on the slide, mirroring the author's cleanup pass across the deck.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($i = $s.Comments.Count; $i -ge 1; $i--) {
        $c = $s.Comments.Item($i)

        for ($j = $c.Replies.Count; $j -ge 1; $j--) {
            $c.Replies.Item($j).Delete()
        }

        $c.Delete()
    }
}
